$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A width: 16 (bestFit) -> 33.5 (explicit custom width) ---
$ws.Columns.Item(1).ColumnWidth = 32.67

# --- B8: clear the red "duplicate" highlight fill (style 2 -> style 3) ---
$ws.Range("B8").ClearFormats()

# --- Row 26: remove the extra "Eighth" / duplicate group sample row ---
# Clear out the data columns entirely (D26:I26 had Participant Count / flags)
$ws.Range("D26:I26").ClearContents()
# A26 / B26 lose their text values but keep being present (still cleared of any fill)
$ws.Range("A26").Value = ""
$ws.Range("B26").Value = ""
$ws.Range("B26").ClearFormats()

# --- Move the active selection to B8 (matches the new sheet selection) ---
[void]$ws.Range("B8").Select()
